$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '51.107.94'
$ws.Range("E2").Value = '  +0.39%  '
$ws.Range("D3").Value = '2.957.65'
$ws.Range("E3").Value = '  +1.09%  '
$ws.Range("E4").Value = '  +0.09%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '379.57'
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = '  +1.48%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '102.05'
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = '  +1.17%  '
$ws.Range("E7").Value = '  +1.82%  '
$ws.Range("E8").Value = '  +0.01%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.592'
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = '  +1.88%  '
$ws.Range("E10").Value = '  +1.09%  '
$ws.Range("E11").Value = '  -1.12%  '
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.0860'
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = '  +1.83%  '
$ws.Range("E13").Value = '  +6.09%  '
$ws.Range("D14").Value = '3.422.98'
$ws.Range("E14").Value = '  +0.97%  '
$ws.Range("E15").Value = '  +2.22%  '
$ws.Range("D16").Value = '2.960.48'
$ws.Range("E16").Value = '  +0.65%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '11.18'
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = '  -4.67%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.995'
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = '  +2.59%  '
$ws.Range("D19").Value = '51.173.67'
$ws.Range("E19").Value = '  +0.62%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '3.14'
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = '  +0.60%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '12.37'
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = '  -0.34%  '
$ws.Range("E22").Value = '  +0.84%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '70.31'
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = '  +2.80%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '266.93'
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = '  +1.26%  '
$ws.Range("E25").Value = '  +1.75%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '7.82'
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = '  -2.27%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '7.28'
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = '  -3.29%  '
$ws.Range("E28").Value = '  +0.01%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '25.86'
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = '  +1.49%  '
$ws.Range("E30").Value = '  -1.54%  '
$ws.Range("E31").Value = '  -0.30%  '
$ws.Range("E32").Value = '  +3.37%  '
$ws.Range("B33").Value = 'InjectiveProtocol'
$ws.Range("C33").Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '34.41'
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = '  +5.06%  '
$ws.Range("B34").Value = 'OKB'
$ws.Range("C34").Value = 'https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb'
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '51.17'
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = '  +1.22%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '2.05'
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = '  +1.69%  '
$ws.Range("E36").Value = '  -0.82%  '
$ws.Range("E37").Value = '  +0.10%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.24'
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = '  +4.15%  '
$ws.Range("E39").Value = '  +1.11%  '
$ws.Range("E40").Value = '  +3.07%  '
$ws.Range("E41").Value = '  +1.79%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '125.19'
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = '  +4.01%  '
$ws.Range("E43").Value = '  +1.07%  '
$ws.Range("E44").Value = '  +7.85%  '
$ws.Range("E45").Value = '  +2.31%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.273'
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = '  -0.66%  '
$ws.Range("E47").Value = '  -0.24%  '
$ws.Range("E48").Value = '  +3.15%  '
$ws.Range("D49").Value = '2.037.90'
$ws.Range("E49").Value = '  +2.59%  '
$ws.Range("E50").Value = '  -3.08%  '
$ws.Range("E51").Value = '  +6.70%  '
